$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 168
$ws.Range("F5").Value = 3938
$ws.Range("F7").Value = 2597
$ws.Range("F9").Value = 3212
$ws.Range("F10").Value = 535
$ws.Range("F11").Value = 2342
$ws.Range("F15").Value = 470
$ws.Range("F18").Value = 228
$ws.Range("F19").Value = 358
$ws.Range("F20").Value = 314
$ws.Range("F21").Value = 455
$ws.Range("F22").Value = 682
$ws.Range("F23").Value = 1433
$ws.Range("F24").Value = 169
$ws.Range("F27").Value = 145
$ws.Range("F28").Value = 162
$ws.Range("F29").Value = 7
$ws.Range("F30").Value = 76
$ws.Range("F31").Value = 4440
$ws.Range("F32").Value = 4317
$ws.Range("F33").Value = 91
$ws.Range("F34").Value = 280
$ws.Range("F35").Value = 76
$ws.Range("F37").Value = 1164
$ws.Range("F38").Value = 163
$ws.Range("F39").Value = 12
$ws.Range("F40").Value = 499
$ws.Range("F44").Value = 137
$ws.Range("F45").Value = 115
$ws.Range("F46").Value = 45

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 9

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 1036
$ws.Range("F4").Value = 2339
$ws.Range("F5").Value = 31

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1036
$ws.Range("F5").Value = 168
$ws.Range("F6").Value = 3938
$ws.Range("F7").Value = 2597
$ws.Range("F9").Value = 3212
$ws.Range("F11").Value = 535
$ws.Range("F12").Value = 2342
$ws.Range("F16").Value = 470
$ws.Range("F18").Value = 228
$ws.Range("F19").Value = 9
$ws.Range("F20").Value = 358
$ws.Range("F21").Value = 682
$ws.Range("F22").Value = 1433
$ws.Range("F24").Value = 145
$ws.Range("F26").Value = 76
$ws.Range("F29").Value = 4440
$ws.Range("F30").Value = 4317
$ws.Range("F31").Value = 91
$ws.Range("F33").Value = 1164
$ws.Range("F34").Value = 163
$ws.Range("F35").Value = 12
$ws.Range("F38").Value = 499
$ws.Range("F45").Value = 115
$ws.Range("F46").Value = 45
